$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Status" column (header in I1, boolean data in I2) is being dropped from
# the product import template. Deleting the whole column shifts "Home Flag"
# and "Hot Flag" one place to the left (J->I, K->J) and keeps their data intact.
$ws.Range("I1:I2").EntireColumn.Delete()

# Shrink "Table1" so its range/column list no longer spans the removed column.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:J23"))

# Resizing the table keeps only the first 10 column definitions, so make sure
# the last two columns of the table are labelled correctly again.
$ws.Range("I1").Value = "Home Flag"
$ws.Range("J1").Value = "Hot Flag"

# Leave the selection where it ended up at the end of the edit.
$ws.Range("C26").Select() | Out-Null
